$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary row updates ---
$ws.Range("B3").Value = 45
$ws.Range("C3").Value = 4381
$ws.Range("D3").Value = 98.8

$ws.Range("B4").Value = 45
$ws.Range("C4").Value = 4381

# --- Good Drivers table: a new driver (21.40.1.3) is now the top entry,
#     pushing the existing rows (12-17) down to (13-18). ---
$ws.Rows.Item(12).Insert()

# Copy formatting from the row that used to be row 12 (now row 13) onto
# the freshly inserted row 12 so styles match the rest of the table.
$ws.Range("A13:E13").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New top entry for the Good Drivers table.
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").Value = 11128
$ws.Range("D12").Value = 100
$ws.Range("E12").ClearContents()

# Updated client counts for the rows that shifted down.
$ws.Range("B13").Value = 486214
$ws.Range("D13").Value = 99.9

$ws.Range("B14").Value = 79953
$ws.Range("D14").Value = 99.9

$ws.Range("B15").Value = 35355

$ws.Range("B16").Value = 65425

$ws.Range("B17").Value = 117653
